$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $text)
    $ws.Range("ZZ1").Formula = "=" + ('"' + $text.Replace('"','""') + '"')
    $ws.Range("ZZ1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
}

$ws.Range("D2").Value = "25.827.84"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.737.16"
$ws.Range("E3").Value = "  -1.03%  "
Set-TextValue "D4" "0.9995"
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue "D5" "231.69"
$ws.Range("E5").Value = "  -1.98%  "
Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  +0.10%  "
Set-TextValue "D7" "0.5192"
$ws.Range("E7").Value = "  +0.21%  "
Set-TextValue "D8" "0.2765"
$ws.Range("E8").Value = "  +2.67%  "
Set-TextValue "D9" "39.27"
$ws.Range("E9").Value = "  -2.78%  "
Set-TextValue "D10" "0.06133"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").Value = "1.736.54"
$ws.Range("E11").Value = "  -1.04%  "
Set-TextValue "D12" "0.07041"
$ws.Range("E12").Value = "  +0.58%  "
Set-TextValue "D13" "15.07"
$ws.Range("E13").Value = "  -2.41%  "
Set-TextValue "D14" "0.6429"
$ws.Range("E14").Value = "  -0.17%  "
Set-TextValue "D15" "4.524"
$ws.Range("E15").Value = "  +0.80%  "
Set-TextValue "D16" "76.73"
$ws.Range("E16").Value = "  -1.76%  "
Set-TextValue "D17" "0.9974"
$ws.Range("E17").Value = "  -0.23%  "
Set-TextValue "D18" "0.9992"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "25.830.43"
$ws.Range("E19").Value = "  -0.44%  "
Set-TextValue "D20" "11.48"
$ws.Range("E20").Value = "  -1.49%  "
Set-TextValue "D21" "0.000006633"
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("D22").Value = "1.956.05"
$ws.Range("E22").Value = "  -1.28%  "
Set-TextValue "D23" "4.177"
$ws.Range("E23").Value = "  +2.46%  "
Set-TextValue "D24" "8.735"
$ws.Range("E24").Value = "  +4.20%  "
Set-TextValue "D25" "5.171"
$ws.Range("E25").Value = "  -0.30%  "
Set-TextValue "D26" "139.65"
$ws.Range("E26").Value = "  +2.44%  "
Set-TextValue "D27" "1.505"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("E28").Value = "  -0.59%  "
Set-TextValue "D29" "1.784"
$ws.Range("E29").Value = "  -2.12%  "
Set-TextValue "D30" "102.18"
$ws.Range("E30").Value = "  -1.03%  "
Set-TextValue "D31" "0.08301"
$ws.Range("E31").Value = "  -0.67%  "
Set-TextValue "D32" "3.682"
$ws.Range("E32").Value = "  -0.43%  "
Set-TextValue "D33" "3.494"
$ws.Range("E33").Value = "  +2.62%  "
Set-TextValue "D34" "0.04495"
$ws.Range("E34").Value = "  +2.29%  "
Set-TextValue "D35" "2.616"
$ws.Range("E35").Value = "  -1.04%  "
Set-TextValue "D36" "0.9781"
$ws.Range("E36").Value = "  -1.89%  "
Set-TextValue "D37" "0.6175"
$ws.Range("E37").Value = "  +1.79%  "
Set-TextValue "D38" "2.674"
$ws.Range("E38").Value = "  -1.77%  "
Set-TextValue "D39" "0.01576"
$ws.Range("E39").Value = "  +0.93%  "
Set-TextValue "D40" "1.930"
$ws.Range("E40").Value = "  -0.54%  "
Set-TextValue "D41" "0.9979"
$ws.Range("E41").Value = "  -0.23%  "
Set-TextValue "D42" "99.87"
$ws.Range("E42").Value = "  -2.47%  "
Set-TextValue "D43" "0.3850"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("E44").Value = "  +1.36%  "
Set-TextValue "D45" "0.7206"
$ws.Range("E45").Value = "  -3.89%  "
Set-TextValue "D46" "0.05351"
$ws.Range("E46").Value = "  -2.63%  "
Set-TextValue "D47" "0.1128"
$ws.Range("E47").Value = "  +1.57%  "
Set-TextValue "D48" "6.207"
$ws.Range("E48").Value = "  +2.49%  "
Set-TextValue "D49" "53.07"
$ws.Range("E49").Value = "  +0.66%  "
Set-TextValue "D50" "29.95"
$ws.Range("E50").Value = "  -0.83%  "
Set-TextValue "D51" "7.623"
$ws.Range("E51").Value = "  +2.87%  "

$ws.Range("ZZ1").ClearContents()
$excel.CutCopyMode = $false
